# "letzte Anpassung im Übungsszenariendokument"
# - Rename sheet "Mitm" -> "MITM"
# - Fix broken reference in D24 of that sheet: was "=#REF!", should be "=F8"
#   (F8 holds the "Int2 = 40" label used elsewhere as the "previous step" caption)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Mitm")
$ws.Name = "MITM"

$ws.Range("D24").Formula = "=F8"
